$wb = $excel.ActiveWorkbook

# Duplicate the "Swiss" sheet as a template for the new "Portugal" sheet,
# placing it after the last sheet, then rename it.
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy($null, $swiss)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Portugal"

# Update the market name and ticket reference for Portugal.
$ws.Range("B2").Value = "Portugal Market"
$ws.Range("B4").Value = "NGC-3479/T2436"

# Resize the columns to the Portugal layout.
$ws.Columns.Item(1).ColumnWidth = 24.385416666666668
$ws.Columns.Item(2).ColumnWidth = 14.721354166666666
$ws.Columns.Item(3).ColumnWidth = 18.944010416666668
$ws.Columns.Item(4).ColumnWidth = 21.166666666666668

# Expand rows 3-5 to fit wrapped text at the new column widths.
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 28.8

# Reset row 14 back to the sheet's default (un-custom) height.
$ws.Rows.Item(14).AutoFit()

# Selection / active cell on the new sheet.
$ws.Range("B4").Select()

# Update selection on the Czech sheet (was a full-column selection).
$czech = $wb.Worksheets.Item("Czech")
$czech.Activate()
$czech.Range("A1:D14").Select()

# Re-activate Portugal so it is the active tab on save.
$ws.Activate()
